$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2044.3572
$ws.Range("I28").Value = 1569.6666
$ws.Range("K28").Value = 1569.6666
$ws.Range("M28").Value = -1084.6666

$ws.Range("H48").Value = 1148.3334
$ws.Range("I48").Value = 722.5
$ws.Range("K48").Value = 2167.5
$ws.Range("M48").Value = -1875.5

$ws.Range("H56").Value = 1148.3334
$ws.Range("I56").Value = 722.5
$ws.Range("K56").Value = 2167.5
$ws.Range("M56").Value = -1633.5

$ws.Range("H86").Value = 129171520
$ws.Range("I86").Value = 247224180
$ws.Range("J86").Value = 11118861
$ws.Range("K86").Value = 247224180
$ws.Range("L86").Value = 11118861
$ws.Range("M86").Value = -247223057
$ws.Range("N86").Value = -11121107

$ws.Range("H89").Value = 129171520
$ws.Range("I89").Value = 247224180
$ws.Range("J89").Value = 11118861
$ws.Range("K89").Value = 1236120900
$ws.Range("L89").Value = 55594305
$ws.Range("M89").Value = -1236115284
$ws.Range("N89").Value = -55605537

$ws.Range("H113").Value = 93756760
$ws.Range("I113").Value = 2721.8
$ws.Range("J113").Value = 136372220
$ws.Range("K113").Value = 2721.8
$ws.Range("L113").Value = 136372220
$ws.Range("M113").Value = 532.1999999999998
$ws.Range("N113").Value = -136378728

$ws.Range("H116").Value = 35724716
$ws.Range("I116").Value = 250000000
$ws.Range("K116").Value = 250000000
$ws.Range("M116").Value = -249996558

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H132").Value = 1637.0476
$ws.Range("I132").Value = 1503.8462
$ws.Range("K132").Value = 4511.5386
$ws.Range("M132").Value = -1981.5386

$ws.Range("H137").Value = 3081.805
$ws.Range("I137").Value = 2990
$ws.Range("J137").Value = 3258.8572
$ws.Range("K137").Value = 8970
$ws.Range("L137").Value = 9776.571599999999
$ws.Range("M137").Value = -6420
$ws.Range("N137").Value = -14876.5716

$ws.Range("H138").Value = 5956.8066
$ws.Range("I138").Value = 985.3158
$ws.Range("K138").Value = 2955.9474
$ws.Range("M138").Value = 2184.0526

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5265.89
$ws.Range("I32").Value = 4934.948
$ws.Range("J32").Value = 15966.333
$ws.Range("K32").Value = 4934.948
$ws.Range("L32").Value = 15966.333
$ws.Range("M32").Value = -4647.948
$ws.Range("N32").Value = -16540.333

$ws.Range("H61").Value = 6506.5483
$ws.Range("I61").Value = 1748.2354
$ws.Range("J61").Value = 12284.5
$ws.Range("K61").Value = 1748.2354
$ws.Range("L61").Value = 12284.5
$ws.Range("M61").Value = -1536.2354
$ws.Range("N61").Value = -12708.5

$ws.Range("H74").Value = 2346.5173
$ws.Range("I74").Value = 1358.8235
$ws.Range("J74").Value = 3745.75
$ws.Range("K74").Value = 1358.8235
$ws.Range("L74").Value = 3745.75
$ws.Range("M74").Value = -484.8235
$ws.Range("N74").Value = -5493.75

$ws.Range("H77").Value = 2346.5173
$ws.Range("I77").Value = 1358.8235
$ws.Range("J77").Value = 3745.75
$ws.Range("K77").Value = 6794.1175
$ws.Range("L77").Value = 18728.75
$ws.Range("M77").Value = -2426.1175
$ws.Range("N77").Value = -27464.75

$ws.Range("H102").Value = 1159.6666
$ws.Range("I102").Value = 1051.7
$ws.Range("K102").Value = 1051.7
$ws.Range("M102").Value = 570.3

$ws.Range("H122").Value = 3035.5945
$ws.Range("I122").Value = 2025.8462
$ws.Range("J122").Value = 5422.273
$ws.Range("K122").Value = 6077.5386
$ws.Range("L122").Value = 16266.819
$ws.Range("M122").Value = -3627.5386
$ws.Range("N122").Value = -21166.819

$ws.Range("H132").Value = 4795.7827
$ws.Range("I132").Value = 2042.3636
$ws.Range("J132").Value = 7319.75
$ws.Range("K132").Value = 6127.0908
$ws.Range("L132").Value = 21959.25
$ws.Range("M132").Value = -3597.0908
$ws.Range("N132").Value = -27019.25

$ws.Range("H136").Value = 6506.5483
$ws.Range("I136").Value = 1748.2354
$ws.Range("J136").Value = 12284.5
$ws.Range("K136").Value = 5244.706200000001
$ws.Range("L136").Value = 36853.5
$ws.Range("M136").Value = -2694.706200000001
$ws.Range("N136").Value = -41953.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5556964.5
$ws.Range("I20").Value = 7937799
$ws.Range("J20").Value = 1683.4445
$ws.Range("K20").Value = 7937799
$ws.Range("L20").Value = 1683.4445
$ws.Range("M20").Value = -7937552
$ws.Range("N20").Value = -2177.4445

$ws.Range("H64").Value = 11112631
$ws.Range("J64").Value = 1985.2941
$ws.Range("L64").Value = 1985.2941
$ws.Range("N64").Value = -2435.2941

$ws.Range("H67").Value = 11112631
$ws.Range("J67").Value = 1985.2941
$ws.Range("L67").Value = 1985.2941
$ws.Range("N67").Value = -3545.2941

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6096.1665
$ws.Range("I16").Value = 4722.625
$ws.Range("K16").Value = 4722.625
$ws.Range("M16").Value = -4435.625

$ws.Range("H31").Value = 6704.719
$ws.Range("I31").Value = 3213.5
$ws.Range("K31").Value = 3213.5
$ws.Range("M31").Value = -2918.5

$ws.Range("H34").Value = 6704.719
$ws.Range("I34").Value = 3213.5
$ws.Range("K34").Value = 3213.5
$ws.Range("M34").Value = -3011.5

$ws.Range("H86").Value = 78130000
$ws.Range("I86").Value = 156250000
$ws.Range("J86").Value = 9999
$ws.Range("K86").Value = 156250000
$ws.Range("L86").Value = 9999
$ws.Range("M86").Value = -156248877
$ws.Range("N86").Value = -12245

$ws.Range("H89").Value = 78130000
$ws.Range("I89").Value = 156250000
$ws.Range("J89").Value = 9999
$ws.Range("K89").Value = 781250000
$ws.Range("L89").Value = 49995
$ws.Range("M89").Value = -781244384
$ws.Range("N89").Value = -61227

$ws.Range("H113").Value = 6096.1665
$ws.Range("I113").Value = 4722.625
$ws.Range("K113").Value = 4722.625
$ws.Range("M113").Value = -2552.625

$ws.Range("H122").Value = 2939.818
$ws.Range("I122").Value = 1473.75
$ws.Range("J122").Value = 3777.5715
$ws.Range("K122").Value = 4421.25
$ws.Range("L122").Value = 11332.7145
$ws.Range("M122").Value = -1971.25
$ws.Range("N122").Value = -16232.7145

$ws.Range("H132").Value = 5150.0703
$ws.Range("I132").Value = 2732.359
$ws.Range("J132").Value = 10388.444
$ws.Range("K132").Value = 8197.076999999999
$ws.Range("L132").Value = 31165.332
$ws.Range("M132").Value = -5667.076999999999
$ws.Range("N132").Value = -36225.33199999999

$ws.Range("H134").Value = 4028.2168
$ws.Range("I134").Value = 1818.9811
$ws.Range("K134").Value = 5456.9433
$ws.Range("M134").Value = -2921.9433

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 55695052
$ws.Range("I4").Value = 6493187
$ws.Range("K4").Value = 19479561
$ws.Range("M4").Value = -19479449

$ws.Range("H86").Value = 1000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 3000
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -5372

$ws.Range("H89").Value = 1000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 9000
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -20856

$ws.Range("H129").Value = 963.1667
$ws.Range("J129").Value = 1791.25
$ws.Range("L129").Value = 5373.75
$ws.Range("N129").Value = -15373.75

$ws.Range("H131").Value = 4222.2
$ws.Range("J131").Value = 4465.25
$ws.Range("L131").Value = 13395.75
$ws.Range("N131").Value = -23475.75

$ws.Range("H132").Value = 9833.223
$ws.Range("J132").Value = 19166.334
$ws.Range("L132").Value = 172497.006
$ws.Range("N132").Value = -177557.006

$ws.Range("H140").Value = 3142.5715
$ws.Range("I140").Value = 1875
$ws.Range("K140").Value = 5625
$ws.Range("M140").Value = -445

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3583850.5
$ws.Range("J122").Value = 3417.625
$ws.Range("L122").Value = 10252.875
$ws.Range("N122").Value = -15152.875

$ws.Range("H132").Value = 8413.895
$ws.Range("I132").Value = 4530.625
$ws.Range("J132").Value = 15070.929
$ws.Range("K132").Value = 13591.875
$ws.Range("L132").Value = 45212.787
$ws.Range("M132").Value = -11061.875
$ws.Range("N132").Value = -50272.787

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7346.6875
$ws.Range("I40").Value = 5448.6
$ws.Range("J40").Value = 8209.454
$ws.Range("K40").Value = 5448.6
$ws.Range("L40").Value = 8209.454
$ws.Range("M40").Value = -5312.6
$ws.Range("N40").Value = -8481.454

$ws.Range("H46").Value = 11114086
$ws.Range("I46").Value = 1325
$ws.Range("J46").Value = 18522592
$ws.Range("K46").Value = 1325
$ws.Range("L46").Value = 18522592
$ws.Range("M46").Value = -1137
$ws.Range("N46").Value = -18522968

$ws.Range("H122").Value = 5818.0645
$ws.Range("I122").Value = 4353.222
$ws.Range("J122").Value = 7846.3076
$ws.Range("K122").Value = 13059.666
$ws.Range("L122").Value = 23538.9228
$ws.Range("M122").Value = -10609.666
$ws.Range("N122").Value = -28438.9228

$ws.Range("H132").Value = 6949505
$ws.Range("I132").Value = 10640755
$ws.Range("K132").Value = 31922265
$ws.Range("M132").Value = -31919735

$ws.Range("H136").Value = 5965.6
$ws.Range("I136").Value = 1557.069
$ws.Range("K136").Value = 4671.207
$ws.Range("M136").Value = -2121.207

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 19608652
$ws.Range("I107").Value = 391.83334
$ws.Range("K107").Value = 1175.50002
$ws.Range("M107").Value = 744.4999800000001

$ws.Range("H122").Value = 113992.195
$ws.Range("I122").Value = 155705.84
$ws.Range("K122").Value = 467117.52
$ws.Range("M122").Value = -464667.52

$ws.Range("H132").Value = 13901352
$ws.Range("I132").Value = 18524174
$ws.Range("K132").Value = 55572522
$ws.Range("M132").Value = -55569992
